$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the fixed "date" placeholder text (7/13/2013 -> 7/14/2013) on the
#    Slide Master and on every Slide Layout (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {}
        if ($isDate) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/13/2013") {
                $tr.Text = "7/14/2013"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 1: rename "TEO" -> "PEO" (TextBox 29).
#    Replace just the leading "T" with "P" so the run is split in two,
#    matching the author's in-place retype of the first letter.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "TEO") {
        $firstChar = $shp.TextFrame.TextRange.Characters(1, 1)
        $firstChar.Text = "P"
    }
}
